# Apply the "Office Theme" colour scheme to the presentation's theme
# (swap away from the "Integral" theme colours it currently has).
#
# PowerPoint's VBA/COM RGB() packs components as R + G*256 + B*65536,
# i.e. little-endian BGR when viewed as a hex integer.
function ToRgbLong($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Index -> (name, target RRGGBB) following the 12-slot a:clrScheme order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$tcs.Colors(1).RGB  = ToRgbLong 0x00 0x00 0x00   # dk1      000000
$tcs.Colors(2).RGB  = ToRgbLong 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Colors(3).RGB  = ToRgbLong 0x44 0x54 0x6A   # dk2      44546A
$tcs.Colors(4).RGB  = ToRgbLong 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Colors(5).RGB  = ToRgbLong 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Colors(6).RGB  = ToRgbLong 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Colors(7).RGB  = ToRgbLong 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Colors(8).RGB  = ToRgbLong 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Colors(9).RGB  = ToRgbLong 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Colors(10).RGB = ToRgbLong 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Colors(11).RGB = ToRgbLong 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Colors(12).RGB = ToRgbLong 0x95 0x4F 0x72   # folHlink 954F72
